$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    $c = $ws.Cells.Item($row, $col)
    $c.Value = "'" + $text
    $c.Style = "Normal"
}

Set-TextCell 2 4 '61.784.10'
Set-TextCell 2 5 '  -0.14%  '
Set-TextCell 3 4 '3.403.91'
Set-TextCell 3 5 '  -0.28%  '
Set-TextCell 4 5 '  -0.10%  '
Set-TextCell 5 4 '408.27'
Set-TextCell 5 5 '  +0.75%  '
Set-TextCell 6 4 '127.89'
Set-TextCell 6 5 '  -3.25%  '
Set-TextCell 7 5 '  +6.81%  '
Set-TextCell 8 5 '  -0.10%  '
Set-TextCell 9 4 '0.727'
Set-TextCell 9 5 '  +6.26%  '
Set-TextCell 10 5 '  +9.71%  '
Set-TextCell 11 4 '42.30'
Set-TextCell 11 5 '  +1.10%  '
Set-TextCell 12 5 '  -0.24%  '
Set-TextCell 13 4 '9.03'
Set-TextCell 13 5 '  +7.18%  '
Set-TextCell 14 4 '3.938.32'
Set-TextCell 14 5 '  -0.52%  '
Set-TextCell 15 4 '21.15'
Set-TextCell 15 5 '  +6.93%  '
Set-TextCell 16 4 '0.0000203'
Set-TextCell 16 5 '  +40.71%  '
Set-TextCell 17 4 '3.370.59'
Set-TextCell 17 5 '  -1.86%  '
Set-TextCell 18 4 '12.04'
Set-TextCell 18 5 '  +3.19%  '
Set-TextCell 19 5 '  +5.44%  '
Set-TextCell 20 4 '61.727.41'
Set-TextCell 20 5 '  -0.25%  '
Set-TextCell 21 4 '454.46'
Set-TextCell 21 5 '  +46.09%  '
Set-TextCell 22 4 '91.42'
Set-TextCell 22 5 '  +9.54%  '
Set-TextCell 23 5 '  -0.63%  '
Set-TextCell 24 4 '12.84'
Set-TextCell 24 5 '  +0.11%  '
Set-TextCell 25 5 '  +3.12%  '
Set-TextCell 26 4 '33.48'
Set-TextCell 26 5 '  +13.16%  '
Set-TextCell 27 4 '8.68'
Set-TextCell 27 5 '  +7.11%  '
Set-TextCell 28 4 '4.76'
Set-TextCell 28 5 '  -0.73%  '
Set-TextCell 29 2 'RenderToken'
Set-TextCell 29 3 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextCell 29 4 '7.53'
Set-TextCell 29 5 '  -2.20%  '
Set-TextCell 30 2 'Toncoin'
Set-TextCell 30 3 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
Set-TextCell 30 4 '2.69'
Set-TextCell 30 5 '  -2.19%  '
Set-TextCell 31 4 '11.91'
Set-TextCell 31 5 '  +5.08%  '
Set-TextCell 32 5 '  -3.58%  '
Set-TextCell 33 4 '42.65'
Set-TextCell 33 5 '  +0.02%  '
Set-TextCell 34 5 '  -1.67%  '
Set-TextCell 35 5 '  -0.04%  '
Set-TextCell 36 4 '0.0493'
Set-TextCell 36 5 '  +2.19%  '
Set-TextCell 37 4 '53.21'
Set-TextCell 37 5 '  +3.95%  '
Set-TextCell 38 4 '0.998'
Set-TextCell 38 5 '  -0.13%  '
Set-TextCell 39 5 '  -0.87%  '
Set-TextCell 40 5 '  +7.23%  '
Set-TextCell 41 4 '2.89'
Set-TextCell 41 5 '  -0.84%  '
Set-TextCell 42 5 '  -3.98%  '
Set-TextCell 43 4 '140.56'
Set-TextCell 43 5 '  +1.43%  '
Set-TextCell 44 4 '4.21'
Set-TextCell 44 5 '  +6.35%  '
Set-TextCell 45 4 '1.97'
Set-TextCell 45 5 '  -0.44%  '
Set-TextCell 46 4 '2.40'
Set-TextCell 46 5 '  +7.91%  '
Set-TextCell 47 4 '16.37'
Set-TextCell 47 5 '  -1.79%  '
Set-TextCell 48 5 '  +4.90%  '
Set-TextCell 49 4 '3.747.26'
Set-TextCell 50 4 '2.103.51'
Set-TextCell 50 5 '  -0.32%  '
Set-TextCell 51 4 '106.58'
Set-TextCell 51 5 '  +28.51%  '
